$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the comma-separated IA Control lists in column A for the affected rows
$ws.Range("A2").Value = 'AU-4 (1),AU-4'
$ws.Range("A3").Value = 'AU-4,AU-14 (1)'
$ws.Range("A4").Value = 'CM-6 b,AU-4'
$ws.Range("A5").Value = 'CM-6 b,SC-5 (2),SC-5'
$ws.Range("A6").Value = 'AU-8 b,AU-7 a,AC-6 (8),AU-7 b,AU-12 (3),CM-5 (1),AC-6 (9)'
$ws.Range("A7").Value = 'AU-12 c,AU-8 b,AU-7 a,CM-6 b,AU-7 b,AU-12 (3),CM-5 (1),AU-12 a'
$ws.Range("A11").Value = 'IA-2 (11),IA-2 (12)'
$ws.Range("A16").Value = 'CM-6 b,CM-7 (2)'
$ws.Range("A21").Value = 'CM-6 b,CM-7 (2)'
$ws.Range("A22").Value = 'CM-6 b,CM-7 (2)'
$ws.Range("A37").Value = 'AC-7 a,AC-7 b'
$ws.Range("A38").Value = 'AC-7 a,AC-7 b'
$ws.Range("A39").Value = 'AC-7 a,AC-7 b'
$ws.Range("A40").Value = 'AC-7 a,AC-7 b'
$ws.Range("A45").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A46").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A47").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A48").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A49").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A50").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A51").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A52").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A53").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A54").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A55").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A56").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A57").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A58").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A59").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A60").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A61").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A62").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A63").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A64").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A65").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A66").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A67").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A68").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A69").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A70").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A71").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A72").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A73").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A74").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A75").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A76").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A77").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A78").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A79").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A80").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A81").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A82").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A83").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A84").Value = 'MA-4 (1) (a),AU-3,AU-3 (1)'
$ws.Range("A85").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A86").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A87").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A88").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A89").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A90").Value = 'MA-4 (1) (a),AU-12 c,AU-3 (1)'
$ws.Range("A91").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A92").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A93").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A94").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A95").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A96").Value = 'MA-4 (1) (a),AU-12 c,AU-3,AU-3 (1)'
$ws.Range("A97").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A98").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A99").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A100").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A101").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A102").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AC-2 (4),AU-12 a'
$ws.Range("A103").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AC-2 (4),AU-12 a'
$ws.Range("A104").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AC-2 (4),AU-12 a'
$ws.Range("A105").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AC-2 (4)'
$ws.Range("A106").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AC-2 (4),AU-12 a'
$ws.Range("A107").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AC-2 (4),AU-12 a'
$ws.Range("A108").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AC-2 (4),AU-12 a'
$ws.Range("A109").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AC-2 (4),AU-12 a'
$ws.Range("A110").Value = 'AU-12 c,AU-3 (1),AU-3,MA-4 (1) (a),AC-2 (4),AU-12 a'
$ws.Range("A111").Value = 'AU-12 c,AU-3 (1),AU-14 (1),AU-3,MA-4 (1) (a),AU-12 a'
$ws.Range("A119").Value = 'MA-4 (1) (a),AU-12 c,AU-3,AU-12 a'
$ws.Range("A120").Value = 'MA-4 (1) (a),AU-12 c,AU-3,AU-12 a'
$ws.Range("A126").Value = 'AC-6 (9),AC-2 (4),CM-5 (1),AU-12 c'
$ws.Range("A128").Value = 'CM-6 b,IA-5 (1) (b),IA-5 (1) (a)'
$ws.Range("A132").Value = 'AC-17 (2),SC-13,MA-4 c,SC-8'
$ws.Range("A133").Value = 'MA-4 (7),AC-12,SC-10,MA-4 e'
$ws.Range("A134").Value = 'AC-12,SC-10'
$ws.Range("A135").Value = 'AC-12,SC-10'
$ws.Range("A137").Value = 'AU-7 a,AU-3 (1),AU-7 (1),AU-14 (1),CM-6 b,AU-3,MA-4 (1) (a),CM-5 (1),AU-12 a,AU-6 (4)'
$ws.Range("A140").Value = 'AU-9 (3),AU-9'
$ws.Range("A141").Value = 'AU-9 (3),AU-9'
$ws.Range("A142").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A143").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A144").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A145").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A146").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A147").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A148").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A149").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A150").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A151").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A152").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A153").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A154").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A165").Value = 'SC-8 (2),SC-8 (1),SC-8'
$ws.Range("A166").Value = 'SC-8 (2),SC-8 (1),SC-8'
$ws.Range("A171").Value = 'AC-11 b,AC-11 a'
$ws.Range("A172").Value = 'AC-11 b,AC-11 a'
$ws.Range("A173").Value = 'AC-11 b,AC-11 a'
$ws.Range("A178").Value = 'AU-4 (1),AU-6 (4),CM-6 b'
$ws.Range("A179").Value = 'AC-17 (1),AC-17 (9),CM-6 b,CM-7 b'
$ws.Range("A180").Value = 'AC-17 (1),CM-6 b,CM-7 b'
$ws.Range("A181").Value = 'SI-11 b,AU-9'
$ws.Range("A182").Value = 'SI-11 b,AU-9'
$ws.Range("A183").Value = 'SI-11 b,AU-9'
$ws.Range("A184").Value = 'SI-11 b,AU-9'
$ws.Range("A185").Value = 'SI-11 b,AU-9'
$ws.Range("A192").Value = 'CM-6 b,AU-3'
$ws.Range("A198").Value = 'AU-4 (1),AU-3'
$ws.Range("A205").Value = 'AU-4 (1),CM-6 b'
$ws.Range("A210").Value = 'AC-8 c 1, AC-8 c 2, AC-8 c 3,AC-8 a'
$ws.Range("A211").Value = 'AC-8 c 1, AC-8 c 2, AC-8 c 3,AC-8 a'
$ws.Range("A212").Value = 'AC-8 c 1, AC-8 c 2, AC-8 c 3,AC-8 a'
$ws.Range("A213").Value = 'AC-8 c 1, AC-8 c 2, AC-8 c 3,AC-8 a'
$ws.Range("A214").Value = 'AC-6 (9),AC-2 (4),AU-12 c'
$ws.Range("A219").Value = 'IA-2 (5),CM-6 b'
$ws.Range("A220").Value = 'IA-2,IA-2 (4),IA-2 (3),IA-2 (5),IA-2 (2)'
$ws.Range("A221").Value = 'IA-2,IA-2 (4),IA-2 (3),IA-2 (5),IA-2 (2)'
$ws.Range("A222").Value = 'AC-18 (1),SC-8 (1),SC-8'
$ws.Range("A229").Value = 'IA-7,CM-7 a'
$ws.Range("A230").Value = 'SC-13,MA-4 (6)'
$ws.Range("A231").Value = 'AC-17 (2),MA-4 (6)'
$ws.Range("A232").Value = 'SC-13,MA-4 (6)'
$ws.Range("A240").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A241").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A243").Value = 'CM-6 b,SI-16,SC-2'
$ws.Range("A262").Value = 'AU-5 (1),AU-5 a'
$ws.Range("A268").Value = 'IA-2 (2),CM-6 b'
$ws.Range("A269").Value = 'IA-2 (1),IA-2 (2),IA-2 (3),IA-2 (4)'
$ws.Range("A271").Value = 'MA-4 (1) (a),AU-12 c'
$ws.Range("A274").Value = 'SC-4,CM-6 b'
$ws.Range("A275").Value = 'SC-4,SC-2'
$ws.Range("A276").Value = 'SC-4,SC-2'
$ws.Range("A279").Value = 'CM-6 b,AU-12 a'
$ws.Range("A282").Value = 'CM-6 b,CM-5 (3)'
$ws.Range("A296").Value = 'IA-2 (11),IA-2 (12)'
$ws.Range("A298").Value = 'IA-2 (1),IA-2 (11),IA-2 (12)'
$ws.Range("A308").Value = 'AU-8 b,AU-8 (1) (a),AU-8 (1) (b)'
$ws.Range("A327").Value = 'CM-5 (1),AU-12 c'
$ws.Range("A329").Value = 'AU-5 b,AU-5 a'
$ws.Range("A344").Value = 'AC-17 (1),CM-7 b'
$ws.Range("A345").Value = 'AC-18 (1),CM-7 a'
$ws.Range("A346").Value = 'CM-6 b,IA-5 (1) (c),CM-7 a'
$ws.Range("A361").Value = 'CM-6 b,CM-7 a'
$ws.Range("A366").Value = 'SI-16,CM-7 a'
$ws.Range("A373").Value = 'CM-6 b,CM-7 a'
$ws.Range("A374").Value = 'CM-6 b,CM-7 a'
$ws.Range("A375").Value = 'CM-6 b,CM-7 a'
$ws.Range("A390").Value = 'CM-6 b,IA-5 (1) (a)'
$ws.Range("A397").Value = 'CM-6 b,SI-16'
$ws.Range("A400").Value = 'CM-6 b,SC-3'
$ws.Range("A401").Value = 'CM-6 b,SC-3'
$ws.Range("A402").Value = 'CM-6 b,SC-3'
$ws.Range("A447").Value = 'CM-6 b,IA-5 (1) (c)'
$ws.Range("A523").Value = 'CM-6 b,SC-2'
$ws.Range("A524").Value = 'CM-6 b,SC-2'
$ws.Range("A540").Value = 'CM-6 b,SI-2 (2)'
$ws.Range("A549").Value = 'CM-6 b,SI-2 (2)'

# Row 168: update Check (K) text and add Fix (M) text
$check168 = @'
To verify that BIND uses the system crypto policy, check out that the BIND config file
 /etc/named.conf  contains the  include "/etc/crypto-policies/back-ends/bind.config"; 
directive:
 $ sudo grep 'include "/etc/crypto-policies/back-ends/bind.config";' /etc/named.conf 
Verify that the directive is at the bottom of the  options  section of the config file.
If BIND is installed and the BIND config file doesn't contain the
 include "/etc/crypto-policies/back-ends/bind.config";  directive then this is a finding.
'@

$fix168 = @'
Configure BIND to use the system crypto policy.
Add the following line to the "options" section in "/etc/named.conf":
include "/etc/crypto-policies/back-ends/bind.config";
'@

$ws.Range("K168").Value = $check168
$ws.Range("M168").Value = $fix168
